$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 6).Value = 3.8
$ws.Cells.Item(3, 7).Value = 5.8
$ws.Cells.Item(3, 9).Value = 2.12
$ws.Cells.Item(3, 10).Value = 3.35
$ws.Cells.Item(3, 11).Value = 5
$ws.Cells.Item(3, 14).Value = 2.06
$ws.Cells.Item(3, 16).Value = 2.06
$ws.Cells.Item(3, 17).Value = 1.75
$ws.Cells.Item(3, 22).Value = 1.89
$ws.Cells.Item(3, 23).Value = 1.2

# Row 4
$ws.Cells.Item(4, 6).Value = 2.04
$ws.Cells.Item(4, 7).Value = 2.74
$ws.Cells.Item(4, 8).Value = 2.74
$ws.Cells.Item(4, 9).Value = 3.9
$ws.Cells.Item(4, 10).Value = 3.15
$ws.Cells.Item(4, 11).Value = 7
$ws.Cells.Item(4, 14).Value = 1.94
$ws.Cells.Item(4, 15).Value = 1.25
$ws.Cells.Item(4, 16).Value = 1.94
$ws.Cells.Item(4, 17).Value = 1.7
$ws.Cells.Item(4, 18).Value = 1.18
$ws.Cells.Item(4, 19).Value = 1.7
$ws.Cells.Item(4, 20).Value = 1.01
$ws.Cells.Item(4, 21).Value = 1.01
$ws.Cells.Item(4, 22).Value = 1.35
$ws.Cells.Item(4, 23).Value = 1.57
$ws.Cells.Item(4, 24).Value = 1000
$ws.Cells.Item(4, 25).Value = 1000
$ws.Cells.Item(4, 28).Value = 1000
$ws.Cells.Item(4, 29).Value = 1000
$ws.Cells.Item(4, 30).Value = 1000
$ws.Cells.Item(4, 32).Value = 1000
$ws.Cells.Item(4, 33).Value = 1000
$ws.Cells.Item(4, 34).Value = 1000

# Row 5
$ws.Cells.Item(5, 14).Value = 1.3
$ws.Cells.Item(5, 16).Value = 1.3
$ws.Cells.Item(5, 19).Value = 1.42

# Row 6
$ws.Cells.Item(6, 6).Value = 2.26
$ws.Cells.Item(6, 8).Value = 2.88
$ws.Cells.Item(6, 9).Value = 3.8
$ws.Cells.Item(6, 10).Value = 3.4
$ws.Cells.Item(6, 11).Value = 4.5
$ws.Cells.Item(6, 14).Value = 3.5
$ws.Cells.Item(6, 15).Value = 1.31
$ws.Cells.Item(6, 16).Value = 1.9
$ws.Cells.Item(6, 17).Value = 1.91
$ws.Cells.Item(6, 18).Value = 1.33
$ws.Cells.Item(6, 19).Value = 3.1
$ws.Cells.Item(6, 20).Value = 1.61
$ws.Cells.Item(6, 21).Value = 1.97
$ws.Cells.Item(6, 22).Value = 1.38

# Row 7
$ws.Cells.Item(7, 6).Value = 2.38
$ws.Cells.Item(7, 7).Value = 3.05
$ws.Cells.Item(7, 8).Value = 2.58
$ws.Cells.Item(7, 10).Value = 3.2
$ws.Cells.Item(7, 11).Value = 4.7
$ws.Cells.Item(7, 12).Value = 1.01
$ws.Cells.Item(7, 13).Value = 1.04
$ws.Cells.Item(7, 14).Value = 3.85
$ws.Cells.Item(7, 15).Value = 1.24
$ws.Cells.Item(7, 16).Value = 2.14
$ws.Cells.Item(7, 17).Value = 1.72
$ws.Cells.Item(7, 18).Value = 1.45
$ws.Cells.Item(7, 19).Value = 2.54
$ws.Cells.Item(7, 20).Value = 1.51
$ws.Cells.Item(7, 21).Value = 2.14
$ws.Cells.Item(7, 22).Value = 1.44
$ws.Cells.Item(7, 23).Value = 1.51

# Row 8
$ws.Cells.Item(8, 6).Value = 3.25
$ws.Cells.Item(8, 8).Value = 2.2
$ws.Cells.Item(8, 12).Value = 1.36
$ws.Cells.Item(8, 13).Value = 1.07
$ws.Cells.Item(8, 14).Value = 3.15
$ws.Cells.Item(8, 16).Value = 1.74
$ws.Cells.Item(8, 17).Value = 2.08
$ws.Cells.Item(8, 19).Value = 3.75
$ws.Cells.Item(8, 24).Value = 15
$ws.Cells.Item(8, 25).Value = 11
$ws.Cells.Item(8, 26).Value = 17.5
$ws.Cells.Item(8, 27).Value = 40
$ws.Cells.Item(8, 28).Value = 15
$ws.Cells.Item(8, 29).Value = 9.2
$ws.Cells.Item(8, 30).Value = 14
$ws.Cells.Item(8, 31).Value = 34
$ws.Cells.Item(8, 32).Value = 30
$ws.Cells.Item(8, 33).Value = 18.5
$ws.Cells.Item(8, 34).Value = 23
$ws.Cells.Item(8, 35).Value = 55
$ws.Cells.Item(8, 36).Value = 75
$ws.Cells.Item(8, 37).Value = 60
$ws.Cells.Item(8, 38).Value = 70
$ws.Cells.Item(8, 39).Value = 1000
$ws.Cells.Item(8, 40).Value = 65
$ws.Cells.Item(8, 41).Value = 28

# Row 11
$ws.Cells.Item(11, 6).Value = 2.42
$ws.Cells.Item(11, 9).Value = 3.55
$ws.Cells.Item(11, 10).Value = 3.1
$ws.Cells.Item(11, 11).Value = 3.4

# Row 16
$ws.Cells.Item(16, 11).Value = 3.4

# Row 17
$ws.Cells.Item(17, 6).Value = 1.7
$ws.Cells.Item(17, 8).Value = 5.3
$ws.Cells.Item(17, 9).Value = 6
$ws.Cells.Item(17, 10).Value = 3.8
$ws.Cells.Item(17, 11).Value = 4.2
$ws.Cells.Item(17, 17).Value = 1.74

# Row 20
$ws.Cells.Item(20, 7).Value = 3.7
$ws.Cells.Item(20, 8).Value = 2.58

# Row 21
$ws.Cells.Item(21, 6).Value = 2.38
$ws.Cells.Item(21, 16).Value = 1.95
$ws.Cells.Item(21, 17).Value = 1.79

# Row 24
$ws.Cells.Item(24, 6).Value = 3.25
$ws.Cells.Item(24, 9).Value = 2.2
$ws.Cells.Item(24, 10).Value = 2.76

# Row 26
$ws.Cells.Item(26, 7).Value = 4.2
$ws.Cells.Item(26, 8).Value = 2.24
$ws.Cells.Item(26, 10).Value = 3
$ws.Cells.Item(26, 11).Value = 3.4
$ws.Cells.Item(26, 16).Value = 1.47
$ws.Cells.Item(26, 17).Value = 2.46

# Row 27
$ws.Cells.Item(27, 6).Value = 3.3
$ws.Cells.Item(27, 7).Value = 3.65
$ws.Cells.Item(27, 8).Value = 2.44
$ws.Cells.Item(27, 9).Value = 2.66
$ws.Cells.Item(27, 16).Value = 1.52
$ws.Cells.Item(27, 17).Value = 2.34

# Row 28
$ws.Cells.Item(28, 6).Value = 2.34
$ws.Cells.Item(28, 7).Value = 2.48
$ws.Cells.Item(28, 8).Value = 3.45
$ws.Cells.Item(28, 9).Value = 3.8
$ws.Cells.Item(28, 10).Value = 3.1
$ws.Cells.Item(28, 11).Value = 3.35
$ws.Cells.Item(28, 16).Value = 1.51

# Row 29
$ws.Cells.Item(29, 7).Value = 2.3
$ws.Cells.Item(29, 17).Value = 2.3
$ws.Cells.Item(29, 20).Value = 1.94
$ws.Cells.Item(29, 21).Value = 1.89
$ws.Cells.Item(29, 24).Value = 13
$ws.Cells.Item(29, 26).Value = 30
$ws.Cells.Item(29, 29).Value = 8.8
$ws.Cells.Item(29, 30).Value = 21
$ws.Cells.Item(29, 31).Value = 70
$ws.Cells.Item(29, 32).Value = 16
$ws.Cells.Item(29, 35).Value = 80
$ws.Cells.Item(29, 36).Value = 980
$ws.Cells.Item(29, 38).Value = 60
$ws.Cells.Item(29, 40).Value = 29
$ws.Cells.Item(29, 41).Value = 95

# Row 30
$ws.Cells.Item(30, 6).Value = 1.77
$ws.Cells.Item(30, 8).Value = 5.9
$ws.Cells.Item(30, 9).Value = 7.4
$ws.Cells.Item(30, 16).Value = 1.51
$ws.Cells.Item(30, 17).Value = 2.68

# Row 31
$ws.Cells.Item(31, 7).Value = 4.8
$ws.Cells.Item(31, 9).Value = 2.12
$ws.Cells.Item(31, 10).Value = 3.25
$ws.Cells.Item(31, 11).Value = 3.55
$ws.Cells.Item(31, 16).Value = 1.66

# Row 32
$ws.Cells.Item(32, 6).Value = 2.36
$ws.Cells.Item(32, 7).Value = 2.64
$ws.Cells.Item(32, 8).Value = 3.25
$ws.Cells.Item(32, 9).Value = 3.85
$ws.Cells.Item(32, 10).Value = 3
$ws.Cells.Item(32, 11).Value = 3.45
$ws.Cells.Item(32, 16).Value = 1.5
$ws.Cells.Item(32, 17).Value = 2.42

# Row 33
$ws.Cells.Item(33, 6).Value = 2.6
$ws.Cells.Item(33, 7).Value = 2.66
$ws.Cells.Item(33, 8).Value = 3.15
$ws.Cells.Item(33, 9).Value = 3.4
$ws.Cells.Item(33, 10).Value = 3
